$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" (columns A:R) ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a new row for client "AREVALO PEÑA JORGE LUIS" under asesor
# "OFICINA-CATAECSA", pushing every row from the old row 245 onward down by one.
$ws1.Rows.Item(245).Insert()

$ws1.Cells.Item(245, 1).Value = "OFICINA-CATAECSA"
$ws1.Cells.Item(245, 2).Value = "AREVALO PEÑA JORGE LUIS"
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item(245, $c).Value = 0
}

# Fix the "N de 275" footer counters (now 276 rows of data) on the totals row,
# which moved from row 277 to row 278 after the insert.
for ($c = 3; $c -le 18; $c++) {
    $cell = $ws1.Cells.Item(278, $c)
    $cell.Value = $cell.Value2.Replace("275", "276")
}

# --- Sheet "VENTA MENSUAL" (columns A:G) ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(245).Insert()

$ws2.Cells.Item(245, 1).Value = "OFICINA-CATAECSA"
$ws2.Cells.Item(245, 2).Value = "AREVALO PEÑA JORGE LUIS"
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item(245, $c).Value = 0
}
